$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 29, shifting existing rows 29:87 down to 30:88.
$ws.Rows(29).Insert()

# Populate the new row 29 with a duplicate of the (now shifted) row 30's
# data -- i.e. the original row-29 record -- except for a new date value.
$ws.Cells.Item(29, 1).Value = 9
$ws.Cells.Item(29, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(29, 3).Value = "Metropolitana"
$ws.Cells.Item(29, 4).Value = 45177
$ws.Cells.Item(29, 5).Value = 13
$ws.Cells.Item(29, 6).Value = 100112010
$ws.Cells.Item(29, 7).Value = "Achicoria"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 70
$ws.Cells.Item(29, 11).Value = 7000
$ws.Cells.Item(29, 12).Value = 7000
$ws.Cells.Item(29, 13).Value = 7000
$ws.Cells.Item(29, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(29, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(29, 16).Value = 438
$ws.Cells.Item(29, 17).Value = 16
$ws.Cells.Item(29, 18).Value = "Hortaliza"
